$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 6676156
$ws.Range("C4").Value = 38837
$ws.Range("D4").Value = 3950354
$ws.Range("E4").Value = 2527678
$ws.Range("G4").Value = 703
$ws.Range("H4").Value = 198124

# --- Row 5: India ---
$ws.Range("D5").Value = 3699298
$ws.Range("E5").Value = 973876

# --- Row 29: Canada ---
$ws.Range("B29").Value = 136141
$ws.Range("C29").Value = 515
$ws.Range("D29").Value = 120075
$ws.Range("E29").Value = 6896

# --- Row 79: Libia ---
$ws.Range("D79").Value = 12100
$ws.Range("E79").Value = 9894

# --- Row 166: Vietnam ---
$ws.Range("D166").Value = 910
$ws.Range("E166").Value = 115

# --- Rows 192-194: Curazao overtakes Brunei and Seychelles in ranking ---
# New order (by total cases, descending): Curazao, Brunei, Seychelles
$ws.Range("A192").Value = "Curazao"
$ws.Range("B192").Value = 145
$ws.Range("C192").Value = 10
$ws.Range("D192").Value = 56
$ws.Range("E192").Value = 88
$ws.Range("H192").Value = 1

$ws.Range("A193").Value = "Brunei"
$ws.Range("B193").Value = 145
$ws.Range("C193").Value = 0
$ws.Range("D193").Value = 139
$ws.Range("H193").Value = 3

$ws.Range("A194").Value = "Seychelles"
$ws.Range("B194").Value = 139
$ws.Range("C194").Value = 1
$ws.Range("D194").Value = 136
$ws.Range("E194").Value = 3
$ws.Range("H194").Value = 0

# --- Row 197: Islas Virgenes Britanicas ---
$ws.Range("B197").Value = 66
$ws.Range("C197").Value = 2
$ws.Range("D197").Value = 37
$ws.Range("E197").Value = 28

# --- Update "last refreshed" timestamp string (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Septiembre de 2020 a las 03:10"
